# Add a new "2022-Q3" quarterly sheet (right after the "总计" summary sheet)
# and record its totals on the "总计" sheet as a new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet, positioned right after "总计".
#    Cloning the "2022-Q2" sheet (same 8-column layout/formatting) and
#    then overwriting its values keeps every bit of formatting (borders,
#    bold header, page margins, outline props, ...) intact.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet    = $wb.Worksheets.Item("2022-Q2")

$q2Sheet.Copy($null, $totalSheet) | Out-Null
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# The source sheet had 8 data rows; 2022-Q3 only has 7, so drop the extra.
$newSheet.Rows.Item(9).Delete() | Out-Null

# Header row (unchanged text, but re-assert it explicitly for clarity).
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holdings for 2022-Q3 (columns B,D-G are text, matching the source data)
$rows = @(
    @(0, "501079", "大成科创主题混合（LOF）A",   "10.00", "85.01", "9.59", "0.9590", 2),
    @(1, "010695", "华夏磐益一年定期开放混合",     "15.90", "99.95", "4.17", "0.6630", 9),
    @(2, "012473", "大成成长回报六个月持有混合A", "7.48",  "80.86", "7.84", "0.5864", 3),
    @(3, "010371", "大成成长进取混合A",           "3.61",  "80.75", "7.12", "0.2570", 3),
    @(4, "010372", "大成成长进取混合C",           "1.52",  "80.75", "7.12", "0.1082", 3),
    @(5, "012474", "大成成长回报六个月持有混合C", "0.37",  "80.86", "7.84", "0.0290", 3),
    @(6, "016198", "大成科创主题混合（LOF）C",   "0.01",  "85.01", "9.59", "0.0010", 2)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    # Column B (fund code, e.g. "010695") is numeric-looking text whose
    # leading zeros must survive — force text storage (without leaving a
    # lingering custom style behind on the cell).
    $codeCell = $newSheet.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[1]
    $codeCell.Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = $row[2]

    # Columns D/E/F/G hold numeric-looking text ("10.00", "0.9590", ...) —
    # force text storage (so trailing zeros survive) without leaving a
    # lingering custom style behind on the cell.
    $textRange = $newSheet.Range($newSheet.Cells.Item($r, 4), $newSheet.Cells.Item($r, 7))
    $textRange.NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $textRange.Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q3 and
#    bump the existing rows' serial numbers down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert() | Out-Null
$totalSheet.Range("B2:D2").ClearFormats() | Out-Null

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 2.6

# Give the new serial-number cell (A2) the same style as the other
# serial-number cells in column A.
$totalSheet.Range("A3").Copy() | Out-Null
$totalSheet.Range("A2").PasteSpecial(-4122) | Out-Null
$totalSheet.Application.CutCopyMode = $false

# The rows that used to be 2/3/4 are now 3/4/5 - their serial numbers
# (column A) each increment by one.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
